$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard Boti")

$ws.Range("D14").Value = 0.6127811275912152
$ws.Range("D14").NumberFormat = "0.00%"

$ws.Range("D15").Value = 2.342400814482517
$ws.Range("D15").NumberFormat = "0.00"

$ws.Range("D16").Value = 0.6788680632120544
$ws.Range("D16").NumberFormat = "0.00%"

$ws.Range("D17").NumberFormat = "0.00%"
